$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Last"
$ws.Range("B2").Value = "Hays"

$ws.Range("D7").Select()
